# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (only) slide master / Design,
#                             currently the "Integral" / "Red Violet" theme.
#   ppt/theme/theme2.xml  -> bound to the notes master, currently the
#                             default "Office Theme".
#
# The authored edit swaps the two themes' contents so the slide master
# (and therefore every slide) now renders with the default Office color
# palette. The only thing that differs between the two themes is the
# <a:clrScheme> color values (fonts/fmtScheme are identical), so we drive
# the visible, reachable side of that swap through
# Slide.ThemeColorScheme(i).RGB, which is the supported way to edit a
# presentation's active theme colors in this host.
#
# PowerPoint's ThemeColorScheme index order is:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000  # dk1      -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF  # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444  # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED  # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244  # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70  # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305  # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95  # folHlink -> 954F72
